# Aula 1 / Etapa 1 - "Correção dos Títulos das Aulas"
# Slide 2 (index 2 in Slides collection) is the cover slide for this stage;
# its background band, title textbox and subtitle textbox are repositioned
# and the title text is corrected.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)

# --- Shape 5 (id 214): big dark background rectangle - tiny 1 EMU reflow ---
$shBg = $s.Shapes.Item(5)
$shBg.Top = 4.511811023622047

# --- Shape 9 (id 218): "Aula 1| Etapa 1: <title>" textbox ---
$shTitle = $s.Shapes.Item(9)
$shTitle.Left = 36.81496242992126
$shTitle.Top = 94.77141962283466
$shTitle.Width = 670.9133858267717
$shTitle.Height = 155.68134308267716

$tf = $shTitle.TextFrame
$tr = $tf.TextRange
$titlePara = $tr.Paragraphs(2)

# Replace "Apresentação" (12 chars) with "Introdução " (11 chars)
$run1 = $titlePara.Characters(1, 12)
$run1.Text = "Introdução "

# Replace " do curso" (9 chars, now starting at position 12) with the new
# second half of the title
$run2 = $titlePara.Characters(12, 9)
$run2.Text = "ao Curso e Conceitos Básicos"

# --- Shape 10 (id 219): "Apresentação" subtitle textbox ---
$shSub = $s.Shapes.Item(10)
$shSub.Left = 36.673228346456696
$shSub.Top = 253.39874015748032
